# Rename the two "SourceKey" header labels to "BusinessKey" on Sheet1,
# per the commit "moved staging files StagingTemplates directory":
#   CustomReportSourceKey -> CustomReportBusinessKey
#   IndicatorSourceKey    -> IndicatorBusinessKey
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "CustomReportBusinessKey"
$ws.Range("C2").Value = "IndicatorBusinessKey"
